# "Fixing Process - 1a"
#
# Insert a new "Relief Amount" column before the existing "Special
# Compensation" column (column G) on Sheet1, and populate it with the
# relief amounts for the two data rows. Everything that used to live in
# columns G..O shifts one column to the right (H..P) automatically when
# the column is inserted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at G; this shifts old G:O -> H:P, preserving their
# values/styles/widths.
$ws.Columns.Item(7).Insert()

# New column header + data for "Relief Amount".
$ws.Range("G1").Value = "Relief Amount"
$ws.Range("G2").Value = 1000
$ws.Range("G3").Value = 2000

# Give the new column (and the now-unformatted "Reason" column to its
# left) explicit widths, matching the template's layout.
$ws.Columns.Item(6).ColumnWidth = 12.666666666666666
$ws.Columns.Item(7).ColumnWidth = 15.5

# Match the saved selection/active cell from the edited workbook.
$ws.Range("G12").Select()
